# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# The document embeds a sample Java exception stack trace. The JDK
# internal frames (java.base/java.util.concurrent.ThreadPoolExecutor...
# and java.base/java.lang.Thread.run) reference line numbers that moved
# after the Apache POI upgrade (4.1.0 -> 5.2.3) regenerated the fixture.
# Update the three line numbers in place.

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2

$replacements = @(
    @{ Old = "ThreadPoolExecutor.java:1130"; New = "ThreadPoolExecutor.java:1136" },
    @{ Old = "ThreadPoolExecutor.java:630";  New = "ThreadPoolExecutor.java:635"  },
    @{ Old = "Thread.java:832";              New = "Thread.java:833"             }
)

foreach ($rep in $replacements) {
    $found = $d.Content.Find.Execute(
        $rep.Old, $true, $false, $false, $false, $false,
        $true, $wdFindContinue, $false, $rep.New, $wdReplaceAll)
    if (-not $found) {
        throw "Could not find text to replace: $($rep.Old)"
    }
}
